$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Wipe the old long-format table (A1:D7); this keeps the existing ---
# --- per-cell formatting in place (row1 = bold/centered header style, ---
# --- rows2-7 = vertical-centered wrap style) so most cells need no   ---
# --- further format work.                                           ---
$ws.Range("A1:D7").ClearContents()

# --- Row 2 becomes the new header row -> give it the bold/centered   ---
# --- header style (copy the format straight from A1, which already  ---
# --- carries it, rather than toggling properties one-by-one so we   ---
# --- don't fork off any stray intermediate cell styles).             ---
$ws.Range("A1").Copy()
$ws.Range("A2:F2").PasteSpecial(-4122)   # xlPasteFormats

# --- The two brand-new columns (E & F) need the same "data" style    ---
# --- used by the rest of the table body, copied from D3.             ---
$ws.Range("D3").Copy()
$ws.Range("E3:F7").PasteSpecial(-4122)   # xlPasteFormats

$ws.Application.CutCopyMode = $false

# --- Row 2: header labels ---
$ws.Range("A2").Value = "Name"
$ws.Range("B2").Value = "Math"
$ws.Range("C2").Value = "English"
$ws.Range("D2").Value = "Science"
$ws.Range("E2").Value = "History"
$ws.Range("F2").Value = "Computer"

# --- Column A: student names (written top to bottom) ---
$ws.Range("A3").Value = "John"
$ws.Range("A4").Value = "Mary"
$ws.Range("A5").Value = "Alice"
$ws.Range("A6").Value = "Bob"
$ws.Range("A7").Value = "David"

# --- Marks for each student ---
$ws.Range("B3").Value = 85
$ws.Range("C3").Value = 78
$ws.Range("D3").Value = 92
$ws.Range("E3").Value = 65
$ws.Range("F3").Value = 88

$ws.Range("B4").Value = 35
$ws.Range("C4").Value = 42
$ws.Range("D4").Value = 38
$ws.Range("E4").Value = 28
$ws.Range("F4").Value = 45

$ws.Range("B5").Value = 100
$ws.Range("C5").Value = 95
$ws.Range("D5").Value = 90
$ws.Range("E5").Value = 98
$ws.Range("F5").Value = 99

$ws.Range("B6").Value = 41
$ws.Range("C6").Value = 33
$ws.Range("D6").Value = 40
$ws.Range("E6").Value = 50
$ws.Range("F6").Value = 30

$ws.Range("B7").Value = 60
$ws.Range("C7").Value = 55
$ws.Range("D7").Value = 70
$ws.Range("E7").Value = 65
$ws.Range("F7").Value = 80

# --- Row heights: drop the old custom (wrapped-text) heights back to default ---
$ws.Rows("1:7").AutoFit()

# --- Column F needs an explicit custom width now that it holds data ---
$ws.Columns("F").ColumnWidth = 8.75

# --- Selection ends up parked on K11, matching the saved view ---
$ws.Range("K11").Select()
